# Regenerate orders with updated distance/size codes.
# Pure text substitution across every string-bearing cell in the sheet:
#   D64 -> D69, D80 -> D86, D51 -> D55 (Distance codes)
#   S30 -> S31                        (Size code; S20/S25 unchanged)
# These tokens are non-overlapping substrings (e.g. "Face06_D64_S25",
# "Fixation_D64_l.png", "D64"), so a straightforward Range.Replace pass
# over the used range reproduces every changed cell/shared-string exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$rng.Replace("D64", "D69")
$rng.Replace("D80", "D86")
$rng.Replace("D51", "D55")
$rng.Replace("S30", "S31")
